$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force a numeric-looking string (e.g. "576.54") to be stored as
    # text rather than being coerced into a Double by the Value setter.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "66.914.33"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.074.71"
$ws.Range("E3").Value = "  -1.06%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws.Range("D5") "576.54"
$ws.Range("E5").Value = "  -0.16%  "
Set-TextValue $ws.Range("D6") "167.68"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.071.65"
$ws.Range("E8").Value = "  -1.06%  "
Set-TextValue $ws.Range("D10") "6.39"
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  -1.26%  "
Set-TextValue $ws.Range("D12") "0.471"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("E13").Value = "  -1.72%  "
Set-TextValue $ws.Range("D14") "35.97"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("E15").Value = "  -1.89%  "
Set-TextValue $ws.Range("D18") "7.02"
$ws.Range("E18").Value = "  -1.36%  "
Set-TextValue $ws.Range("D19") "16.79"
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("D20").Value = "3.068.11"
$ws.Range("E20").Value = "  -1.30%  "
Set-TextValue $ws.Range("D21") "485.75"
$ws.Range("E21").Value = "  +1.91%  "
Set-TextValue $ws.Range("D22") "0.688"
$ws.Range("E22").Value = "  -3.40%  "
Set-TextValue $ws.Range("D23") "7.70"
$ws.Range("E23").Value = "  -2.89%  "
Set-TextValue $ws.Range("D24") "82.62"
$ws.Range("E24").Value = "  -1.46%  "
Set-TextValue $ws.Range("D25") "12.81"
$ws.Range("E25").Value = "  -4.74%  "
Set-TextValue $ws.Range("D26") "2.22"
$ws.Range("E26").Value = "  -2.82%  "
Set-TextValue $ws.Range("D27") "10.25"
$ws.Range("E27").Value = "  +2.94%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("E30").Value = "  -5.91%  "
Set-TextValue $ws.Range("D31") "2.62"
$ws.Range("E31").Value = "  -1.65%  "
Set-TextValue $ws.Range("D32") "27.70"
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").Value = "0.0₃0904"
$ws.Range("E34").Value = "  -3.25%  "
Set-TextValue $ws.Range("D35") "0.999"
$ws.Range("E35").Value = "  -0.12%  "
Set-TextValue $ws.Range("D36") "5.65"
$ws.Range("E36").Value = "  -3.21%  "
Set-TextValue $ws.Range("D37") "0.954"
$ws.Range("E37").Value = "  -2.18%  "
Set-TextValue $ws.Range("D38") "45.92"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  -4.47%  "
Set-TextValue $ws.Range("D41") "0.303"
$ws.Range("E41").Value = "  -1.98%  "
Set-TextValue $ws.Range("D42") "8.31"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").Value = "2.765.36"
$ws.Range("E43").Value = "  -1.29%  "
Set-TextValue $ws.Range("D44") "369.48"
$ws.Range("E44").Value = "  -1.92%  "
Set-TextValue $ws.Range("D45") "136.04"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  -3.05%  "
Set-TextValue $ws.Range("D47") "2.48"
$ws.Range("E47").Value = "  -2.60%  "
Set-TextValue $ws.Range("D49") "24.58"
$ws.Range("E49").Value = "  -0.85%  "
Set-TextValue $ws.Range("D50") "2.16"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -1.76%  "

# Rows 16 and 17 swap (WrappedBTC <-> WrappedliquidstakedEther2.0)
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.585.87"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.836.38"
$ws.Range("E17").Value = "  +0.19%  "
